# The "Prix Spot" sheet gained a new date column "09-dec" right before the
# "01-oct." block (column EL), pushing that whole October block (and
# everything after it, through the old FP) one column to the right
# (to FQ). We reproduce this by inserting a whole column at EL (column 142)
# and filling the freshly inserted column with the new header (row 1) and
# "-" placeholders (rows 2-25), matching the rest of that date block's style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newColIndex = 142  # column "EL"

# Shift "EL:FP" (and everything in between) one column to the right by
# inserting a new blank column at EL.
$ws.Columns.Item($newColIndex).Insert()

# Populate the newly inserted column.
$ws.Cells.Item(1, $newColIndex).Value = "09-dec"

for ($row = 2; $row -le 25; $row++) {
    $ws.Cells.Item($row, $newColIndex).Value = "-"
}
